$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 value (Prisma_Excel_File for pop1) to point at the new template file
$ws.Range("C2").Value = "\Testdata\Templates\PRISMA\Test_Sachin\12. PRISMA_Pfizer_IC AML Mylotarg.xlsx"

# Remove the pop2 (ICER) block of rows 7-10 entirely
$ws.Rows("7:10").Delete()

# Update the active selection to match the saved view state
$ws.Range("E11").Select()
